$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.618.61'
$ws.Range('E2').Value = '  +3.50%  '
$ws.Range('D3').Value = '3.434.98'
$ws.Range('E3').Value = '  +2.72%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.22'
$ws.Range('E5').Value = '  +2.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '183.65'
$ws.Range('E6').Value = '  +5.18%  '
$ws.Range('E7').Value = '  +1.98%  '
$ws.Range('D8').Value = '3.431.33'
$ws.Range('E8').Value = '  +2.81%  '
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('E10').Value = '  +2.94%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.642'
$ws.Range('E11').Value = '  +1.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '56.10'
$ws.Range('E12').Value = '  +4.38%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000277'
$ws.Range('E13').Value = '  +0.98%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.38'
$ws.Range('E14').Value = '  +3.53%  '
$ws.Range('D15').Value = '3.979.87'
$ws.Range('E15').Value = '  +2.58%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.55'
$ws.Range('E16').Value = '  +2.34%  '
$ws.Range('D17').Value = '3.425.38'
$ws.Range('E17').Value = '  +2.19%  '
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('D19').Value = '66.534.00'
$ws.Range('E19').Value = '  +2.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.04'
$ws.Range('E20').Value = '  +2.79%  '
$ws.Range('E21').Value = '  +2.64%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '492.51'
$ws.Range('E22').Value = '  +8.90%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '16.61'
$ws.Range('E23').Value = '  +15.98%  '
$ws.Range('E24').Value = '  +1.76%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.23'
$ws.Range('E25').Value = '  +3.22%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '89.27'
$ws.Range('E26').Value = '  +2.47%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.06'
$ws.Range('E27').Value = '  +3.03%  '
$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.95'
$ws.Range('E28').Value = '  +2.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.16'
$ws.Range('E29').Value = '  +5.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.33'
$ws.Range('E30').Value = '  +0.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.16'
$ws.Range('E31').Value = '  +8.85%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.67'
$ws.Range('E32').Value = '  +1.83%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '593.07'
$ws.Range('E33').Value = '  +4.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.27'
$ws.Range('E34').Value = '  +3.82%  '
$ws.Range('E35').Value = '  +3.64%  '
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('E37').Value = '  +4.64%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.58'
$ws.Range('E38').Value = '  -1.50%  '
$ws.Range('B39').Value = 'PEPE'
$ws.Range('C39').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D39').Value = '0.0₃0768'
$ws.Range('E39').Value = '  +4.28%  '
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.384'
$ws.Range('E40').Value = '  +4.06%  '
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.07'
$ws.Range('E41').Value = '  +2.16%  '
$ws.Range('D42').Value = '3.163.19'
$ws.Range('E42').Value = '  +3.39%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.93'
$ws.Range('E43').Value = '  +5.09%  '
$ws.Range('E44').Value = '  +2.68%  '
$ws.Range('E45').Value = '  +4.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.80'
$ws.Range('E46').Value = '  +21.64%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.135'
$ws.Range('E47').Value = '  +0.72%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.18'
$ws.Range('E48').Value = '  -0.97%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.75'
$ws.Range('E49').Value = '  +7.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.999'
$ws.Range('E50').Value = '  -0.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '140.59'
$ws.Range('E51').Value = '  +0.06%  '
